# Insert a new data row at row 302 (pushes the existing rows 302:325 down to 303:326)
# and populate it with the new record's values. This mirrors the author's edit of
# adding one more weekly price observation to the "Arándano (blue)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(302).Insert()

$newRow = 302

$ws.Cells.Item($newRow, 1).Value = 9
$ws.Cells.Item($newRow, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value = [DateTime]"2023-03-28"
$ws.Cells.Item($newRow, 5).Value = 13
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100101001
$ws.Cells.Item($newRow, 10).Value = "Arándano (blue)"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 380
$ws.Cells.Item($newRow, 14).Value = 3800
$ws.Cells.Item($newRow, 15).Value = 4000
$ws.Cells.Item($newRow, 16).Value = 3905
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($newRow, 19).Value = 1952
$ws.Cells.Item($newRow, 20).Value = 2

$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
